$wb = $excel.ActiveWorkbook

# --- Phase1 sheet ---
$ws1 = $wb.Worksheets.Item("Phase1")
$ws1.Range("B2").Value = 0.1
$ws1.Range("B3").Value = 0.1
$ws1.Range("B4").Value = 0.05
$ws1.Range("B5").Value = 0.05

# --- Phase2 sheet ---
$ws2 = $wb.Worksheets.Item("Phase2")
$ws2.Range("C2").Value = 0.1
$ws2.Range("C3").Value = 0.1
$ws2.Range("C4").Value = 0.05
$ws2.Range("C5").Value = 0.05

# --- Phase5 sheet ---
$ws5 = $wb.Worksheets.Item("Phase5")

# Update selections on Phase2 and Phase5 first (without activating them),
# then finally activate Phase1 and select B2:B5 so Phase1 ends up as the
# active/selected tab.
$ws2.Activate()
$ws2.Range("C2:C5").Select()

$ws5.Activate()
$ws5.Range("K14").Select()

$ws1.Activate()
$ws1.Range("B2:B5").Select()

$wb.Save()
